$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring column K's formatting in line with column J (new 2021 year column)
$ws.Range("J3:J10").Copy() | Out-Null
$ws.Range("K3:K10").PasteSpecial(-4122) | Out-Null

# Fill in the new column's values
$ws.Range("K4").Value = 2021
$ws.Range("K5").Value = 375
$ws.Range("K6").Value = "-"
$ws.Range("K7").Value = 5
$ws.Range("K8").Value = "-"
$ws.Range("K9").Value = 18
$ws.Range("K10").Value = 150

# Match the selection left behind by the edit (active cell K7)
$ws.Range("K7").Select() | Out-Null
